# Update NATMI Bsg-Slc16a1 LR-pair table with newly computed TPM-based values.
#
# The workbook holds one row per (Sending cluster, Target cluster) pair for the
# Bsg (ligand) -> Slc16a1 (receptor) interaction, for the three clusters
# ECs / FAPs / MuSCs (9 rows total, rows 2-10).
#
# Columns:
#   G = Ligand average expression value          (depends on Sending cluster)
#   H = Ligand total expression value             (depends on Sending cluster)
#   I = Ligand specificity of average expression  (G / sum(G) over sending clusters)
#   J = Ligand specificity of total expression    (H / sum(H) over sending clusters)
#   M = Receptor average expression value         (depends on Target cluster)
#   N = Receptor total expression value           (depends on Target cluster)
#   O = Receptor specificity of average expression (M / sum(M) over target clusters)
#   P = Receptor specificity of total expression   (N / sum(N) over target clusters)
#   Q = Edge average expression weight             (G * M)
#   R = Edge total expression weight               (H * N)
#   S = Edge average expression specificity        (Q / sum(Q) over all 9 pairs)
#   T = Edge total expression specificity          (R / sum(R) over all 9 pairs)
#
# With the new TPM recomputation, only the ligand (Bsg) values for the "ECs"
# sending cluster and the receptor (Slc16a1) values for the "ECs" target
# cluster change; everything derived from them (specificities / edge weights)
# is recalculated accordingly for all 9 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2..10

# Map each data row to its Sending cluster / Target cluster label, read from
# columns A (sending) and D (target).
$sending = @{}
$target = @{}
foreach ($r in $rows) {
    $sending[$r] = $ws.Range("A$r").Value2
    $target[$r] = $ws.Range("D$r").Value2
}

# New ligand average/total expression values, keyed by sending cluster.
# Only "ECs" changes with the refreshed TPM data; FAPs/MuSCs keep their
# existing values.
$newG = @{
    "ECs"   = 33.84781966666667
    "FAPs"  = 74.45592499999999
    "MuSCs" = 32.60069533333333
}
$newH = @{
    "ECs"   = 101.543459
    "FAPs"  = 223.367775
    "MuSCs" = 97.802086
}

# New receptor average/total expression values, keyed by target cluster.
# Only "ECs" changes; FAPs/MuSCs keep their existing values.
$newM = @{
    "ECs"   = 3.393633333333334
    "FAPs"  = 4.323660333333334
    "MuSCs" = 4.152602333333333
}
$newN = @{
    "ECs"   = 10.1809
    "FAPs"  = 12.970981
    "MuSCs" = 12.457807
}

# Write the refreshed base ligand/receptor values into G, H, M, N.
foreach ($r in $rows) {
    $ws.Range("G$r").Value = $newG[$sending[$r]]
    $ws.Range("H$r").Value = $newH[$sending[$r]]
    $ws.Range("M$r").Value = $newM[$target[$r]]
    $ws.Range("N$r").Value = $newN[$target[$r]]
}

# Recompute the ligand specificity values (I, J): value divided by the sum of
# that same value across the three sending clusters.
$sumG = 0.0
$sumH = 0.0
foreach ($c in @("ECs", "FAPs", "MuSCs")) {
    $sumG += $newG[$c]
    $sumH += $newH[$c]
}

# Recompute the receptor specificity values (O, P): value divided by the sum
# of that same value across the three target clusters.
$sumM = 0.0
$sumN = 0.0
foreach ($c in @("ECs", "FAPs", "MuSCs")) {
    $sumM += $newM[$c]
    $sumN += $newN[$c]
}

foreach ($r in $rows) {
    $ws.Range("I$r").Value = $newG[$sending[$r]] / $sumG
    $ws.Range("J$r").Value = $newH[$sending[$r]] / $sumH
    $ws.Range("O$r").Value = $newM[$target[$r]] / $sumM
    $ws.Range("P$r").Value = $newN[$target[$r]] / $sumN
}

# Recompute the edge expression weights (Q, R) as the product of the ligand
# and receptor values for each row, then their derived specificities (S, T)
# as that product divided by the sum of all 9 row products.
$sumQ = 0.0
$sumR = 0.0
foreach ($r in $rows) {
    $q = $newG[$sending[$r]] * $newM[$target[$r]]
    $rr = $newH[$sending[$r]] * $newN[$target[$r]]
    $sumQ += $q
    $sumR += $rr
}

foreach ($r in $rows) {
    $q = $newG[$sending[$r]] * $newM[$target[$r]]
    $rr = $newH[$sending[$r]] * $newN[$target[$r]]
    $ws.Range("Q$r").Value = $q
    $ws.Range("R$r").Value = $rr
    $ws.Range("S$r").Value = $q / $sumQ
    $ws.Range("T$r").Value = $rr / $sumR
}
